$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -11.03379999999999
$ws.Range("A3").Value = -21.22700000000002
$ws.Range("C5").Value = -13.98470000000001
$ws.Range("A14").Value = -20.45289999999998
$ws.Range("A16").Value = -20.35459999999999
$ws.Range("C16").Value = -11.9287
$ws.Range("A21").Value = -21.2591
$ws.Range("A23").Value = -21.24720000000002
$ws.Range("A25").Value = -22.50840000000003
